# Add a new "cluster_class" column (H) to the harmonic-centrality results
# sheet, matching the header style of the existing columns and filling in
# the per-row cluster-class values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H1 -------------------------------------------------
# Copy the formatting (bold, centered, bordered) from the neighboring
# header cell G1, then set the new header text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "cluster_class"

# --- Data values for H2:H63 -----------------------------------------
$clusterClass = @{
    2  = 0
    3  = 0
    4  = -1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = -1
    20 = 0
    21 = -1
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 0
    60 = 0
    61 = 0
    62 = 0
    63 = -0.05
}

foreach ($row in 2..63) {
    $ws.Cells.Item($row, 8).Value = $clusterClass[$row]
}
